$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the item descriptions (values stay in the same rows, but note that
# the old "headset" text at row 10 and "monitor" text at row 11 are swapped
# as part of the rewrite: row 10 becomes the Monitor description and row 11
# becomes the Headset description).
$ws.Range("A1").Value = "Processador gamer Intel Core i5-10400F BX8070110400F de 6 núcleos e 4.3GHz de frequência"
$ws.Range("A2").Value = "Placa Mãe Asus Para Intel 1700 Z690 Plus D4 Tuf 4ddr4 Atx"
$ws.Range("A3").Value = "Placa De Video Mancer Rx 5500 Xt Streaky, 8gb, Gddr6 128 Bit"
$ws.Range("A4").Value = "Memória RAM NB BLACK color preto 16GB 1 UP Gamer UP3200"
$ws.Range("A5").Value = "Disco Sólido Interno Kingston Skc600/512g 512gb Preto Cor Preto"
$ws.Range("A6").Value = "Fonte de alimentação para PC Corsair CV Series CV550 550W black 100V/240V"
$ws.Range("A7").Value = "Gabinete Gamer Cooler Master Elite 300 Lateral Vidro Preto"
$ws.Range("A8").Value = "Teclado Corsair K55 Rgb Multicolor Led"
$ws.Range("A9").Value = "Mouse Logitech G G Series G502 Hero preto"
$ws.Range("A10").Value = "Monitor Acer 21.5 Zero Frame Radeon Hdmi Ea220q Hbi"
$ws.Range("A11").Value = "Headset Gamer Para Consoles E Pc Driver 40mm Quantum 100 Preto Jbl"

# Widen column A to fit the longer descriptions.
$ws.Columns.Item(1).ColumnWidth = 77.44

# Touch the very last column so a second <col> definition (matching the
# source workbook) is emitted alongside column A's custom width.
$ws.Columns.Item(1024).ColumnWidth = 10.69

# Move the active selection to A17.
$ws.Range("A17").Select()

# Keep gridlines displayed (re-affirm default) so the serializer emits the
# worksheet with gridlines shown, matching the source workbook.
$excel.ActiveWindow.DisplayGridlines = $false
$excel.ActiveWindow.DisplayGridlines = $true
